$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.892.94'
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").Value = '1.832.56'
$ws.Range("E3").Value = '  -1.96%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '244.89'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6889'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.000'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07689'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.07%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3050'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.67%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.38'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07806'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.04%  '
$ws.Range("D12").Value = '1.835.98'
$ws.Range("E12").Value = '  -2.95%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.084'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '90.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.46%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6815'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.96%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.444'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.94%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008303'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").Value = '28.884.97'
$ws.Range("E18").Value = '  -1.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.42'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.83%  '
$ws.Range("D20").Value = '2.075.17'
$ws.Range("E20").Value = '  -3.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.99%  '
$ws.Range("E22").Value = '  +0.04%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.459'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.000'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.04%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1477'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.58%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.18'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.809'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.19'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.23%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.544'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.53%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.213'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.26%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.152'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.42%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.183'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05113'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7662'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.840'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.10%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.140'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.44%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.696'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01844'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.71%  '
$ws.Range("D39").Value = '1.217.65'
$ws.Range("E39").Value = '  -4.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.697'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.77%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9399'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '108.76'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9998'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.729'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.81%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000123'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.12%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5164'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.62%  '
$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.519'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.04%  '
$ws.Range("D48").Value = '1.975.58'
$ws.Range("E48").Value = '  -3.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.20'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -9.39%  '
$ws.Range("E50").Value = '  -3.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4188'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.47%  '
